# "all multi party house contraints defined"
# Rearrange the Household-flow columns (L:P) on the "Sets" sheet so the
# "to Car" / "Battery to" values move up from column P into L/M, and the
# old L/M values move down into M/P - completing the multi-party-house
# constraint matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sets")
$ws.Activate()

# --- Row 1 (headers): rotate L1 -> M1 -> P1 -> L1 ---
$ws.Range("L1").Value = "to Car"
$ws.Range("M1").Value = "Battery to"
$ws.Range("P1").Value = "HP to"

# --- Row 2 ---
$ws.Range("L2").Value = "PV"
$ws.Range("M2").Value = "Car"
$ws.Range("P2").Value = "Household"

# --- Row 3 ---
$ws.Range("L3").Value = "Electric Grid"
$ws.Range("P3").ClearContents()

# --- Row 4 ---
$ws.Range("L4").Value = "Battery"
$ws.Range("P4").ClearContents()

# Column J ("Electric Grid to") was widened (best-fit) after the edits.
$ws.Columns.Item(10).ColumnWidth = 12

# Restore the view to the top-left and move the active selection to D2.
$ws.Range("A1").Select() | Out-Null
$ws.Range("D2").Select() | Out-Null
